$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9-10 down to 10-11.
$ws.Rows(9).Insert()

# Populate the newly inserted row 9 with the new representative microgrid town.
$ws.Range("A9").Value = "North"
$ws.Range("B9").Value = "Nord-Est"
$ws.Range("C9").Value = 19.516156630000001
$ws.Range("D9").Value = -71.937673329999996
$ws.Range("E9").Value = "Perches"

# Apply wrap-text / vertically centered style (matching the other data rows' look)
# to the new row's populated cells.
$ws.Range("A9:E9").WrapText = $true
$ws.Range("A9:E9").VerticalAlignment = -4108
$ws.Range("A9:E9").Font.ThemeColor = 1

Write-Host "Done"
